# Auto-generated PowerShell Excel COM-interop script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Reorder match results within existing rows (F:V columns) ---
# Row 74 <- source data from former row 75
$ws.Range('F74').Value = 'QPR'
$ws.Range('G74').Value = 1
$ws.Range('H74').Value = 'Swansea'
$ws.Range('I74').Value = 1
$ws.Range('J74').Value = 2.99
$ws.Range('K74').Value = '16/09/2023 17:39'
$ws.Range('L74').Value = 2.76
$ws.Range('M74').Value = '19/09/2023 20:40'
$ws.Range('N74').Value = 3.4
$ws.Range('O74').Value = '16/09/2023 17:39'
$ws.Range('P74').Value = 3.43
$ws.Range('Q74').Value = '19/09/2023 20:38'
$ws.Range('R74').Value = 2.47
$ws.Range('S74').Value = '16/09/2023 17:39'
$ws.Range('T74').Value = 2.68
$ws.Range('U74').Value = '19/09/2023 20:40'
$ws.Range('V74').Value = 'https://www.betexplorer.com/football/england/championship/qpr-swansea/QoErbhG8/'

# Row 75 <- source data from former row 76
$ws.Range('F75').Value = 'Preston'
$ws.Range('G75').Value = 2
$ws.Range('H75').Value = 'Birmingham'
$ws.Range('I75').Value = 1
$ws.Range('J75').Value = 2.24
$ws.Range('K75').Value = '16/09/2023 17:13'
$ws.Range('L75').Value = 2.42
$ws.Range('M75').Value = '19/09/2023 20:41'
$ws.Range('N75').Value = 3.3
$ws.Range('O75').Value = '16/09/2023 17:13'
$ws.Range('P75').Value = 3.09
$ws.Range('Q75').Value = '19/09/2023 20:41'
$ws.Range('R75').Value = 3.53
$ws.Range('S75').Value = '16/09/2023 17:13'
$ws.Range('T75').Value = 3.44
$ws.Range('U75').Value = '19/09/2023 20:41'
$ws.Range('V75').Value = 'https://www.betexplorer.com/football/england/championship/preston-birmingham/YD0BdQMD/'

# Row 76 <- source data from former row 74
$ws.Range('F76').Value = 'Southampton'
$ws.Range('G76').Value = 0
$ws.Range('H76').Value = 'Ipswich'
$ws.Range('I76').Value = 1
$ws.Range('J76').Value = 2.29
$ws.Range('K76').Value = '16/09/2023 17:13'
$ws.Range('L76').Value = 2.2
$ws.Range('M76').Value = '19/09/2023 20:44'
$ws.Range('N76').Value = 3.74
$ws.Range('O76').Value = '16/09/2023 17:13'
$ws.Range('P76').Value = 3.71
$ws.Range('Q76').Value = '19/09/2023 20:44'
$ws.Range('R76').Value = 3.04
$ws.Range('S76').Value = '16/09/2023 17:13'
$ws.Range('T76').Value = 3.28
$ws.Range('U76').Value = '19/09/2023 20:44'
$ws.Range('V76').Value = 'https://www.betexplorer.com/football/england/championship/southampton-ipswich/vDSrvI6n/'

# Row 77 <- source data from former row 78
$ws.Range('F77').Value = 'Bristol City'
$ws.Range('G77').Value = 4
$ws.Range('H77').Value = 'Plymouth'
$ws.Range('I77').Value = 1
$ws.Range('J77').Value = 1.98
$ws.Range('K77').Value = '16/09/2023 17:13'
$ws.Range('L77').Value = 1.62
$ws.Range('M77').Value = '19/09/2023 20:40'
$ws.Range('N77').Value = 3.77
$ws.Range('O77').Value = '16/09/2023 17:13'
$ws.Range('P77').Value = 4.42
$ws.Range('Q77').Value = '19/09/2023 20:41'
$ws.Range('R77').Value = 3.78
$ws.Range('S77').Value = '16/09/2023 17:13'
$ws.Range('T77').Value = 5.34
$ws.Range('U77').Value = '19/09/2023 20:41'
$ws.Range('V77').Value = 'https://www.betexplorer.com/football/england/championship/bristol-city-plymouth/UuzVmYID/'

# Row 78 <- source data from former row 77
$ws.Range('F78').Value = 'Cardiff'
$ws.Range('G78').Value = 3
$ws.Range('H78').Value = 'Coventry'
$ws.Range('I78').Value = 2
$ws.Range('J78').Value = 2.67
$ws.Range('K78').Value = '16/09/2023 20:12'
$ws.Range('L78').Value = 2.77
$ws.Range('M78').Value = '19/09/2023 20:44'
$ws.Range('N78').Value = 3.37
$ws.Range('O78').Value = '16/09/2023 20:12'
$ws.Range('P78').Value = 3.39
$ws.Range('Q78').Value = '19/09/2023 20:15'
$ws.Range('R78').Value = 2.76
$ws.Range('S78').Value = '16/09/2023 20:12'
$ws.Range('T78').Value = 2.69
$ws.Range('U78').Value = '19/09/2023 20:44'
$ws.Range('V78').Value = 'https://www.betexplorer.com/football/england/championship/cardiff-coventry/neFvaY02/'

# Row 87 <- source data from former row 93
$ws.Range('F87').Value = 'Ipswich'
$ws.Range('G87').Value = 4
$ws.Range('H87').Value = 'Blackburn'
$ws.Range('I87').Value = 3
$ws.Range('J87').Value = 1.72
$ws.Range('K87').Value = '20/09/2023 00:43'
$ws.Range('L87').Value = 1.66
$ws.Range('M87').Value = '23/09/2023 15:56'
$ws.Range('N87').Value = 4.11
$ws.Range('O87').Value = '20/09/2023 00:43'
$ws.Range('P87').Value = 4.38
$ws.Range('Q87').Value = '23/09/2023 15:56'
$ws.Range('R87').Value = 4.69
$ws.Range('S87').Value = '20/09/2023 00:43'
$ws.Range('T87').Value = 5.02
$ws.Range('U87').Value = '23/09/2023 15:59'
$ws.Range('V87').Value = 'https://www.betexplorer.com/football/england/championship/ipswich-blackburn/dAEEiyyn/'

# Row 88 <- source data from former row 92
$ws.Range('F88').Value = 'Leeds'
$ws.Range('G88').Value = 3
$ws.Range('H88').Value = 'Watford'
$ws.Range('I88').Value = 0
$ws.Range('J88').Value = 1.71
$ws.Range('K88').Value = '20/09/2023 00:43'
$ws.Range('L88').Value = 1.77
$ws.Range('M88').Value = '23/09/2023 15:59'
$ws.Range('N88').Value = 3.96
$ws.Range('O88').Value = '20/09/2023 00:43'
$ws.Range('P88').Value = 3.97
$ws.Range('Q88').Value = '23/09/2023 15:58'
$ws.Range('R88').Value = 4.96
$ws.Range('S88').Value = '20/09/2023 00:43'
$ws.Range('T88').Value = 4.7
$ws.Range('U88').Value = '23/09/2023 15:59'
$ws.Range('V88').Value = 'https://www.betexplorer.com/football/england/championship/leeds-watford/vZ3Jjejh/'

# Row 89 <- source data from former row 94
$ws.Range('F89').Value = 'Leicester'
$ws.Range('G89').Value = 1
$ws.Range('H89').Value = 'Bristol City'
$ws.Range('I89').Value = 0
$ws.Range('J89').Value = 1.63
$ws.Range('K89').Value = '20/09/2023 00:43'
$ws.Range('L89').Value = 1.66
$ws.Range('M89').Value = '23/09/2023 15:58'
$ws.Range('N89').Value = 4.36
$ws.Range('O89').Value = '20/09/2023 00:43'
$ws.Range('P89').Value = 4.28
$ws.Range('Q89').Value = '23/09/2023 15:59'
$ws.Range('R89').Value = 5.19
$ws.Range('S89').Value = '20/09/2023 00:43'
$ws.Range('T89').Value = 5.15
$ws.Range('U89').Value = '23/09/2023 15:59'
$ws.Range('V89').Value = 'https://www.betexplorer.com/football/england/championship/leicester-bristol-city/j97NkF5b/'

# Row 92 <- source data from former row 89
$ws.Range('F92').Value = 'Rotherham'
$ws.Range('G92').Value = 1
$ws.Range('H92').Value = 'Preston'
$ws.Range('I92').Value = 1
$ws.Range('J92').Value = 2.95
$ws.Range('K92').Value = '18/09/2023 19:42'
$ws.Range('L92').Value = 3.52
$ws.Range('M92').Value = '23/09/2023 15:58'
$ws.Range('N92').Value = 3.3
$ws.Range('O92').Value = '18/09/2023 19:42'
$ws.Range('P92').Value = 3.25
$ws.Range('Q92').Value = '23/09/2023 15:50'
$ws.Range('R92').Value = 2.56
$ws.Range('S92').Value = '18/09/2023 19:42'
$ws.Range('T92').Value = 2.29
$ws.Range('U92').Value = '23/09/2023 15:58'
$ws.Range('V92').Value = 'https://www.betexplorer.com/football/england/championship/rotherham-preston/tb0zmDkH/'

# Row 93 <- source data from former row 88
$ws.Range('F93').Value = 'Swansea'
$ws.Range('G93').Value = 3
$ws.Range('H93').Value = 'Sheffield Wed'
$ws.Range('I93').Value = 0
$ws.Range('J93').Value = 1.94
$ws.Range('K93').Value = '20/09/2023 00:43'
$ws.Range('L93').Value = 2.07
$ws.Range('M93').Value = '23/09/2023 15:49'
$ws.Range('N93').Value = 3.62
$ws.Range('O93').Value = '20/09/2023 00:43'
$ws.Range('P93').Value = 3.42
$ws.Range('Q93').Value = '23/09/2023 15:49'
$ws.Range('R93').Value = 4.12
$ws.Range('S93').Value = '20/09/2023 00:43'
$ws.Range('T93').Value = 3.94
$ws.Range('U93').Value = '23/09/2023 15:49'
$ws.Range('V93').Value = 'https://www.betexplorer.com/football/england/championship/swansea-sheffield-wed/IgWkdWpL/'

# Row 94 <- source data from former row 87
$ws.Range('F94').Value = 'West Brom'
$ws.Range('G94').Value = 0
$ws.Range('H94').Value = 'Millwall'
$ws.Range('I94').Value = 0
$ws.Range('J94').Value = 1.98
$ws.Range('K94').Value = '18/09/2023 19:42'
$ws.Range('L94').Value = 2.12
$ws.Range('M94').Value = '23/09/2023 15:59'
$ws.Range('N94').Value = 3.57
$ws.Range('O94').Value = '18/09/2023 19:42'
$ws.Range('P94').Value = 3.43
$ws.Range('Q94').Value = '23/09/2023 15:59'
$ws.Range('R94').Value = 4.08
$ws.Range('S94').Value = '18/09/2023 19:42'
$ws.Range('T94').Value = 3.78
$ws.Range('U94').Value = '23/09/2023 15:59'
$ws.Range('V94').Value = 'https://www.betexplorer.com/football/england/championship/west-brom-millwall/EBbroiKT/'

# Row 100 <- source data from former row 107
$ws.Range('F100').Value = 'Bristol City'
$ws.Range('G100').Value = 2
$ws.Range('H100').Value = 'Stoke'
$ws.Range('I100').Value = 3
$ws.Range('J100').Value = 2.46
$ws.Range('K100').Value = '23/09/2023 17:13'
$ws.Range('L100').Value = 2.69
$ws.Range('M100').Value = '30/09/2023 15:55'
$ws.Range('N100').Value = 3.38
$ws.Range('O100').Value = '23/09/2023 17:13'
$ws.Range('P100').Value = 3.33
$ws.Range('Q100').Value = '30/09/2023 15:51'
$ws.Range('R100').Value = 3.06
$ws.Range('S100').Value = '23/09/2023 17:13'
$ws.Range('T100').Value = 2.81
$ws.Range('U100').Value = '30/09/2023 15:55'
$ws.Range('V100').Value = 'https://www.betexplorer.com/football/england/championship/bristol-city-stoke-city/hxUh5f5o/'

# Row 101 <- source data from former row 106
$ws.Range('F101').Value = 'Cardiff'
$ws.Range('G101').Value = 2
$ws.Range('H101').Value = 'Rotherham'
$ws.Range('I101').Value = 0
$ws.Range('J101').Value = 1.84
$ws.Range('K101').Value = '23/09/2023 17:13'
$ws.Range('L101').Value = 1.59
$ws.Range('M101').Value = '30/09/2023 15:50'
$ws.Range('N101').Value = 3.76
$ws.Range('O101').Value = '23/09/2023 17:13'
$ws.Range('P101').Value = 4.06
$ws.Range('Q101').Value = '30/09/2023 15:51'
$ws.Range('R101').Value = 4.5
$ws.Range('S101').Value = '23/09/2023 17:13'
$ws.Range('T101').Value = 6.43
$ws.Range('U101').Value = '30/09/2023 15:51'
$ws.Range('V101').Value = 'https://www.betexplorer.com/football/england/championship/cardiff-rotherham/0ILfej0R/'

# Row 102 <- source data from former row 104
$ws.Range('F102').Value = 'Huddersfield'
$ws.Range('G102').Value = 1
$ws.Range('H102').Value = 'Ipswich'
$ws.Range('I102').Value = 1
$ws.Range('J102').Value = 3.65
$ws.Range('K102').Value = '24/09/2023 03:13'
$ws.Range('L102').Value = 3.53
$ws.Range('M102').Value = '30/09/2023 15:58'
$ws.Range('N102').Value = 3.82
$ws.Range('O102').Value = '24/09/2023 03:13'
$ws.Range('P102').Value = 3.63
$ws.Range('Q102').Value = '30/09/2023 15:58'
$ws.Range('R102').Value = 2.01
$ws.Range('S102').Value = '24/09/2023 03:13'
$ws.Range('T102').Value = 2.13
$ws.Range('U102').Value = '30/09/2023 15:58'
$ws.Range('V102').Value = 'https://www.betexplorer.com/football/england/championship/huddersfield-ipswich/WSSd4EKi/'

# Row 103 <- source data from former row 105
$ws.Range('F103').Value = 'Hull'
$ws.Range('G103').Value = 1
$ws.Range('H103').Value = 'Plymouth'
$ws.Range('I103').Value = 1
$ws.Range('J103').Value = 1.7
$ws.Range('K103').Value = '24/09/2023 03:13'
$ws.Range('L103').Value = 1.86
$ws.Range('M103').Value = '30/09/2023 15:43'
$ws.Range('N103').Value = 4.23
$ws.Range('O103').Value = '24/09/2023 03:13'
$ws.Range('P103').Value = 3.94
$ws.Range('Q103').Value = '30/09/2023 15:43'
$ws.Range('R103').Value = 4.7
$ws.Range('S103').Value = '24/09/2023 03:13'
$ws.Range('T103').Value = 4.21
$ws.Range('U103').Value = '30/09/2023 15:43'
$ws.Range('V103').Value = 'https://www.betexplorer.com/football/england/championship/hull-city-plymouth/Q3s03YZc/'

# Row 104 <- source data from former row 108
$ws.Range('F104').Value = 'Millwall'
$ws.Range('G104').Value = 0
$ws.Range('H104').Value = 'Swansea'
$ws.Range('I104').Value = 3
$ws.Range('J104').Value = 2.08
$ws.Range('K104').Value = '23/09/2023 17:13'
$ws.Range('L104').Value = 2.31
$ws.Range('M104').Value = '30/09/2023 15:50'
$ws.Range('N104').Value = 3.52
$ws.Range('O104').Value = '23/09/2023 17:13'
$ws.Range('P104').Value = 3.36
$ws.Range('Q104').Value = '30/09/2023 14:18'
$ws.Range('R104').Value = 3.72
$ws.Range('S104').Value = '23/09/2023 17:13'
$ws.Range('T104').Value = 3.36
$ws.Range('U104').Value = '30/09/2023 15:50'
$ws.Range('V104').Value = 'https://www.betexplorer.com/football/england/championship/millwall-swansea/2R7MjA8r/'

# Row 105 <- source data from former row 103
$ws.Range('F105').Value = 'Norwich'
$ws.Range('G105').Value = 2
$ws.Range('H105').Value = 'Birmingham'
$ws.Range('I105').Value = 0
$ws.Range('J105').Value = 1.69
$ws.Range('K105').Value = '23/09/2023 17:13'
$ws.Range('L105').Value = 1.89
$ws.Range('M105').Value = '30/09/2023 15:56'
$ws.Range('N105').Value = 4.11
$ws.Range('O105').Value = '23/09/2023 17:13'
$ws.Range('P105').Value = 3.67
$ws.Range('Q105').Value = '30/09/2023 15:56'
$ws.Range('R105').Value = 4.91
$ws.Range('S105').Value = '23/09/2023 17:13'
$ws.Range('T105').Value = 4.36
$ws.Range('U105').Value = '30/09/2023 15:47'
$ws.Range('V105').Value = 'https://www.betexplorer.com/football/england/championship/norwich-birmingham/4ft42hk4/'

# Row 106 <- source data from former row 102
$ws.Range('F106').Value = 'Preston'
$ws.Range('G106').Value = 0
$ws.Range('H106').Value = 'West Brom'
$ws.Range('I106').Value = 4
$ws.Range('J106').Value = 2.81
$ws.Range('K106').Value = '23/09/2023 17:13'
$ws.Range('L106').Value = 2.65
$ws.Range('M106').Value = '30/09/2023 15:57'
$ws.Range('N106').Value = 3.19
$ws.Range('O106').Value = '23/09/2023 17:13'
$ws.Range('P106').Value = 3.17
$ws.Range('Q106').Value = '30/09/2023 15:37'
$ws.Range('R106').Value = 2.75
$ws.Range('S106').Value = '23/09/2023 17:13'
$ws.Range('T106').Value = 3
$ws.Range('U106').Value = '30/09/2023 15:57'
$ws.Range('V106').Value = 'https://www.betexplorer.com/football/england/championship/preston-west-brom/nJV81C4A/'

# Row 107 <- source data from former row 101
$ws.Range('F107').Value = 'QPR'
$ws.Range('G107').Value = 1
$ws.Range('H107').Value = 'Coventry'
$ws.Range('I107').Value = 3
$ws.Range('J107').Value = 2.94
$ws.Range('K107').Value = '23/09/2023 17:13'
$ws.Range('L107').Value = 2.93
$ws.Range('M107').Value = '30/09/2023 15:56'
$ws.Range('N107').Value = 3.42
$ws.Range('O107').Value = '23/09/2023 17:13'
$ws.Range('P107').Value = 3.34
$ws.Range('Q107').Value = '30/09/2023 15:54'
$ws.Range('R107').Value = 2.52
$ws.Range('S107').Value = '23/09/2023 17:13'
$ws.Range('T107').Value = 2.59
$ws.Range('U107').Value = '30/09/2023 15:56'
$ws.Range('V107').Value = 'https://www.betexplorer.com/football/england/championship/qpr-coventry/0MZC0WJG/'

# Row 108 <- source data from former row 100
$ws.Range('F108').Value = 'Watford'
$ws.Range('G108').Value = 2
$ws.Range('H108').Value = 'Middlesbrough'
$ws.Range('I108').Value = 3
$ws.Range('J108').Value = 2.39
$ws.Range('K108').Value = '23/09/2023 17:13'
$ws.Range('L108').Value = 2.78
$ws.Range('M108').Value = '30/09/2023 15:56'
$ws.Range('N108').Value = 3.55
$ws.Range('O108').Value = '23/09/2023 17:13'
$ws.Range('P108').Value = 3.45
$ws.Range('Q108').Value = '30/09/2023 15:56'
$ws.Range('R108').Value = 2.99
$ws.Range('S108').Value = '23/09/2023 17:13'
$ws.Range('T108').Value = 2.65
$ws.Range('U108').Value = '30/09/2023 15:56'
$ws.Range('V108').Value = 'https://www.betexplorer.com/football/england/championship/watford-middlesbrough/ETFmfUdp/'

# Row 110 <- source data from former row 112
$ws.Range('F110').Value = 'Birmingham'
$ws.Range('G110').Value = 4
$ws.Range('H110').Value = 'Huddersfield'
$ws.Range('I110').Value = 1
$ws.Range('J110').Value = 2.01
$ws.Range('K110').Value = '26/09/2023 20:12'
$ws.Range('L110').Value = 2.17
$ws.Range('M110').Value = '03/10/2023 20:40'
$ws.Range('N110').Value = 3.52
$ws.Range('O110').Value = '26/09/2023 20:12'
$ws.Range('P110').Value = 3.27
$ws.Range('Q110').Value = '03/10/2023 20:40'
$ws.Range('R110').Value = 4.03
$ws.Range('S110').Value = '26/09/2023 20:12'
$ws.Range('T110').Value = 3.83
$ws.Range('U110').Value = '03/10/2023 20:31'
$ws.Range('V110').Value = 'https://www.betexplorer.com/football/england/championship/birmingham-huddersfield/hKEiglCj/'

# Row 111 <- source data from former row 113
$ws.Range('F111').Value = 'Ipswich'
$ws.Range('G111').Value = 3
$ws.Range('H111').Value = 'Hull'
$ws.Range('I111').Value = 0
$ws.Range('J111').Value = 1.71
$ws.Range('K111').Value = '29/09/2023 15:42'
$ws.Range('L111').Value = 1.83
$ws.Range('M111').Value = '03/10/2023 20:08'
$ws.Range('N111').Value = 4.1
$ws.Range('O111').Value = '29/09/2023 15:42'
$ws.Range('P111').Value = 3.85
$ws.Range('Q111').Value = '03/10/2023 20:32'
$ws.Range('R111').Value = 4.8
$ws.Range('S111').Value = '29/09/2023 15:42'
$ws.Range('T111').Value = 4.44
$ws.Range('U111').Value = '03/10/2023 20:32'
$ws.Range('V111').Value = 'https://www.betexplorer.com/football/england/championship/ipswich-hull-city/6N2ZHXli/'

# Row 112 <- source data from former row 110
$ws.Range('F112').Value = 'Middlesbrough'
$ws.Range('G112').Value = 2
$ws.Range('H112').Value = 'Cardiff'
$ws.Range('I112').Value = 0
$ws.Range('J112').Value = 1.64
$ws.Range('K112').Value = '26/09/2023 20:12'
$ws.Range('L112').Value = 1.85
$ws.Range('M112').Value = '03/10/2023 20:30'
$ws.Range('N112').Value = 4.33
$ws.Range('O112').Value = '26/09/2023 20:12'
$ws.Range('P112').Value = 3.72
$ws.Range('Q112').Value = '03/10/2023 20:40'
$ws.Range('R112').Value = 5.11
$ws.Range('S112').Value = '26/09/2023 20:12'
$ws.Range('T112').Value = 4.56
$ws.Range('U112').Value = '03/10/2023 20:44'
$ws.Range('V112').Value = 'https://www.betexplorer.com/football/england/championship/middlesbrough-cardiff/W0CQkUNl/'

# Row 113 <- source data from former row 111
$ws.Range('F113').Value = 'Plymouth'
$ws.Range('G113').Value = 0
$ws.Range('H113').Value = 'Millwall'
$ws.Range('I113').Value = 2
$ws.Range('J113').Value = 2.61
$ws.Range('K113').Value = '29/09/2023 15:42'
$ws.Range('L113').Value = 2.6
$ws.Range('M113').Value = '03/10/2023 20:38'
$ws.Range('N113').Value = 3.52
$ws.Range('O113').Value = '29/09/2023 15:42'
$ws.Range('P113').Value = 3.55
$ws.Range('Q113').Value = '03/10/2023 20:32'
$ws.Range('R113').Value = 2.75
$ws.Range('S113').Value = '29/09/2023 15:42'
$ws.Range('T113').Value = 2.77
$ws.Range('U113').Value = '03/10/2023 20:38'
$ws.Range('V113').Value = 'https://www.betexplorer.com/football/england/championship/plymouth-millwall/Ic6wHi4c/'

# Row 116 <- source data from former row 120
$ws.Range('F116').Value = 'Swansea'
$ws.Range('G116').Value = 2
$ws.Range('H116').Value = 'Norwich'
$ws.Range('I116').Value = 1
$ws.Range('J116').Value = 2.68
$ws.Range('K116').Value = '27/09/2023 20:12'
$ws.Range('L116').Value = 2.82
$ws.Range('M116').Value = '04/10/2023 20:37'
$ws.Range('N116').Value = 3.49
$ws.Range('O116').Value = '27/09/2023 20:12'
$ws.Range('P116').Value = 3.36
$ws.Range('Q116').Value = '04/10/2023 19:57'
$ws.Range('R116').Value = 2.68
$ws.Range('S116').Value = '27/09/2023 20:12'
$ws.Range('T116').Value = 2.66
$ws.Range('U116').Value = '04/10/2023 20:43'
$ws.Range('V116').Value = 'https://www.betexplorer.com/football/england/championship/swansea-norwich/AiBUllwe/'

# Row 117 <- source data from former row 119
$ws.Range('F117').Value = 'Sunderland'
$ws.Range('G117').Value = 2
$ws.Range('H117').Value = 'Watford'
$ws.Range('I117').Value = 0
$ws.Range('J117').Value = 2.2
$ws.Range('K117').Value = '27/09/2023 20:12'
$ws.Range('L117').Value = 2.16
$ws.Range('M117').Value = '04/10/2023 20:43'
$ws.Range('N117').Value = 3.6
$ws.Range('O117').Value = '27/09/2023 20:12'
$ws.Range('P117').Value = 3.46
$ws.Range('Q117').Value = '04/10/2023 20:43'
$ws.Range('R117').Value = 3.34
$ws.Range('S117').Value = '27/09/2023 20:12'
$ws.Range('T117').Value = 3.61
$ws.Range('U117').Value = '04/10/2023 20:43'
$ws.Range('V117').Value = 'https://www.betexplorer.com/football/england/championship/sunderland-watford/0pvV67tj/'

# Row 118 <- source data from former row 116
$ws.Range('F118').Value = 'Coventry'
$ws.Range('G118').Value = 1
$ws.Range('H118').Value = 'Blackburn'
$ws.Range('I118').Value = 0
$ws.Range('J118').Value = 2.03
$ws.Range('K118').Value = '27/09/2023 20:12'
$ws.Range('L118').Value = 2.14
$ws.Range('M118').Value = '04/10/2023 20:37'
$ws.Range('N118').Value = 3.75
$ws.Range('O118').Value = '27/09/2023 20:12'
$ws.Range('P118').Value = 3.71
$ws.Range('Q118').Value = '04/10/2023 20:37'
$ws.Range('R118').Value = 3.62
$ws.Range('S118').Value = '27/09/2023 20:12'
$ws.Range('T118').Value = 3.43
$ws.Range('U118').Value = '04/10/2023 20:44'
$ws.Range('V118').Value = 'https://www.betexplorer.com/football/england/championship/coventry-blackburn/8najEklG/'

# Row 119 <- source data from former row 117
$ws.Range('F119').Value = 'Leeds'
$ws.Range('G119').Value = 1
$ws.Range('H119').Value = 'QPR'
$ws.Range('I119').Value = 0
$ws.Range('J119').Value = 1.32
$ws.Range('K119').Value = '29/09/2023 15:42'
$ws.Range('L119').Value = 1.32
$ws.Range('M119').Value = '04/10/2023 20:22'
$ws.Range('N119').Value = 5.74
$ws.Range('O119').Value = '29/09/2023 15:42'
$ws.Range('P119').Value = 5.85
$ws.Range('Q119').Value = '04/10/2023 20:30'
$ws.Range('R119').Value = 9.67
$ws.Range('S119').Value = '29/09/2023 15:42'
$ws.Range('T119').Value = 9.470000000000001
$ws.Range('U119').Value = '04/10/2023 20:42'
$ws.Range('V119').Value = 'https://www.betexplorer.com/football/england/championship/leeds-qpr/nwbfD93M/'

# Row 120 <- source data from former row 118
$ws.Range('F120').Value = 'Leicester'
$ws.Range('G120').Value = 3
$ws.Range('H120').Value = 'Preston'
$ws.Range('I120').Value = 0
$ws.Range('J120').Value = 1.44
$ws.Range('K120').Value = '29/09/2023 15:42'
$ws.Range('L120').Value = 1.5
$ws.Range('M120').Value = '04/10/2023 20:18'
$ws.Range('N120').Value = 4.69
$ws.Range('O120').Value = '29/09/2023 15:42'
$ws.Range('P120').Value = 4.44
$ws.Range('Q120').Value = '04/10/2023 20:40'
$ws.Range('R120').Value = 7.53
$ws.Range('S120').Value = '29/09/2023 15:42'
$ws.Range('T120').Value = 7.17
$ws.Range('U120').Value = '04/10/2023 20:35'
$ws.Range('V120').Value = 'https://www.betexplorer.com/football/england/championship/leicester-preston/U1maCTIS/'

# Row 126 <- source data from former row 128
$ws.Range('F126').Value = 'Plymouth'
$ws.Range('G126').Value = 1
$ws.Range('H126').Value = 'Swansea'
$ws.Range('I126').Value = 3
$ws.Range('J126').Value = 2.55
$ws.Range('K126').Value = '04/10/2023 04:42'
$ws.Range('L126').Value = 2.54
$ws.Range('M126').Value = '07/10/2023 15:47'
$ws.Range('N126').Value = 3.46
$ws.Range('O126').Value = '04/10/2023 04:42'
$ws.Range('P126').Value = 3.62
$ws.Range('Q126').Value = '07/10/2023 15:06'
$ws.Range('R126').Value = 2.87
$ws.Range('S126').Value = '04/10/2023 04:42'
$ws.Range('T126').Value = 2.8
$ws.Range('U126').Value = '07/10/2023 15:47'
$ws.Range('V126').Value = 'https://www.betexplorer.com/football/england/championship/plymouth-swansea/6J4wmS77/'

# Row 127 <- source data from former row 126
$ws.Range('F127').Value = 'Millwall'
$ws.Range('G127').Value = 2
$ws.Range('H127').Value = 'Hull'
$ws.Range('I127').Value = 2
$ws.Range('J127').Value = 1.88
$ws.Range('K127').Value = '30/09/2023 17:13'
$ws.Range('L127').Value = 2.6
$ws.Range('M127').Value = '07/10/2023 15:49'
$ws.Range('N127').Value = 3.78
$ws.Range('O127').Value = '30/09/2023 17:13'
$ws.Range('P127').Value = 3.28
$ws.Range('Q127').Value = '07/10/2023 15:39'
$ws.Range('R127').Value = 4.16
$ws.Range('S127').Value = '30/09/2023 17:13'
$ws.Range('T127').Value = 2.96
$ws.Range('U127').Value = '07/10/2023 15:49'
$ws.Range('V127').Value = 'https://www.betexplorer.com/football/england/championship/millwall-hull-city/6sajAtmr/'

# Row 128 <- source data from former row 127
$ws.Range('F128').Value = 'Leicester'
$ws.Range('G128').Value = 2
$ws.Range('H128').Value = 'Stoke'
$ws.Range('I128').Value = 0
$ws.Range('J128').Value = 1.55
$ws.Range('K128').Value = '04/10/2023 04:42'
$ws.Range('L128').Value = 1.49
$ws.Range('M128').Value = '07/10/2023 15:36'
$ws.Range('N128').Value = 4.44
$ws.Range('O128').Value = '04/10/2023 04:42'
$ws.Range('P128').Value = 4.63
$ws.Range('Q128').Value = '07/10/2023 15:57'
$ws.Range('R128').Value = 6.17
$ws.Range('S128').Value = '04/10/2023 04:42'
$ws.Range('T128').Value = 7.06
$ws.Range('U128').Value = '07/10/2023 15:55'
$ws.Range('V128').Value = 'https://www.betexplorer.com/football/england/championship/leicester-stoke-city/QqoMR8IF/'

# Row 135 <- source data from former row 138
$ws.Range('F135').Value = 'Middlesbrough'
$ws.Range('G135').Value = 1
$ws.Range('H135').Value = 'Birmingham'
$ws.Range('I135').Value = 0
$ws.Range('J135').Value = 1.83
$ws.Range('K135').Value = '09/10/2023 14:42'
$ws.Range('L135').Value = 1.74
$ws.Range('M135').Value = '21/10/2023 15:37'
$ws.Range('N135').Value = 3.78
$ws.Range('O135').Value = '09/10/2023 14:42'
$ws.Range('P135').Value = 3.93
$ws.Range('Q135').Value = '21/10/2023 15:37'
$ws.Range('R135').Value = 4.42
$ws.Range('S135').Value = '09/10/2023 14:42'
$ws.Range('T135').Value = 4.98
$ws.Range('U135').Value = '21/10/2023 15:37'
$ws.Range('V135').Value = 'https://www.betexplorer.com/football/england/championship/middlesbrough-birmingham/rwi1Su9Q/'

# Row 138 <- source data from former row 135
$ws.Range('F138').Value = 'Swansea'
$ws.Range('G138').Value = 1
$ws.Range('H138').Value = 'Leicester'
$ws.Range('I138').Value = 3
$ws.Range('J138').Value = 4.14
$ws.Range('K138').Value = '09/10/2023 14:42'
$ws.Range('L138').Value = 4.37
$ws.Range('M138').Value = '21/10/2023 15:59'
$ws.Range('N138').Value = 3.91
$ws.Range('O138').Value = '09/10/2023 14:42'
$ws.Range('P138').Value = 3.92
$ws.Range('Q138').Value = '21/10/2023 15:51'
$ws.Range('R138').Value = 1.88
$ws.Range('S138').Value = '09/10/2023 14:42'
$ws.Range('T138').Value = 1.83
$ws.Range('U138').Value = '21/10/2023 15:51'
$ws.Range('V138').Value = 'https://www.betexplorer.com/football/england/championship/swansea-leicester/bDiPAVaE/'

# Row 140 <- source data from former row 141
$ws.Range('F140').Value = 'West Brom'
$ws.Range('G140').Value = 0
$ws.Range('H140').Value = 'Plymouth'
$ws.Range('I140').Value = 0
$ws.Range('J140').Value = 1.68
$ws.Range('K140').Value = '09/10/2023 14:42'
$ws.Range('L140').Value = 1.66
$ws.Range('M140').Value = '21/10/2023 15:37'
$ws.Range('N140').Value = 4.29
$ws.Range('O140').Value = '09/10/2023 14:42'
$ws.Range('P140').Value = 4.22
$ws.Range('Q140').Value = '21/10/2023 15:54'
$ws.Range('R140').Value = 4.91
$ws.Range('S140').Value = '09/10/2023 14:42'
$ws.Range('T140').Value = 5.23
$ws.Range('U140').Value = '21/10/2023 15:54'
$ws.Range('V140').Value = 'https://www.betexplorer.com/football/england/championship/west-brom-plymouth/MiUmJyWC/'

# Row 141 <- source data from former row 144
$ws.Range('F141').Value = 'Bristol City'
$ws.Range('G141').Value = 1
$ws.Range('H141').Value = 'Coventry'
$ws.Range('I141').Value = 0
$ws.Range('J141').Value = 2.22
$ws.Range('K141').Value = '09/10/2023 14:42'
$ws.Range('L141').Value = 2.95
$ws.Range('M141').Value = '21/10/2023 15:47'
$ws.Range('N141').Value = 3.59
$ws.Range('O141').Value = '09/10/2023 14:42'
$ws.Range('P141').Value = 3.39
$ws.Range('Q141').Value = '21/10/2023 15:47'
$ws.Range('R141').Value = 3.28
$ws.Range('S141').Value = '09/10/2023 14:42'
$ws.Range('T141').Value = 2.54
$ws.Range('U141').Value = '21/10/2023 15:29'
$ws.Range('V141').Value = 'https://www.betexplorer.com/football/england/championship/bristol-city-coventry/Y79A5I1D/'

# Row 144 <- source data from former row 140
$ws.Range('F144').Value = 'Huddersfield'
$ws.Range('G144').Value = 2
$ws.Range('H144').Value = 'QPR'
$ws.Range('I144').Value = 1
$ws.Range('J144').Value = 1.93
$ws.Range('K144').Value = '09/10/2023 14:42'
$ws.Range('L144').Value = 1.82
$ws.Range('M144').Value = '21/10/2023 15:39'
$ws.Range('N144').Value = 3.64
$ws.Range('O144').Value = '09/10/2023 14:42'
$ws.Range('P144').Value = 3.77
$ws.Range('Q144').Value = '21/10/2023 15:48'
$ws.Range('R144').Value = 4.15
$ws.Range('S144').Value = '09/10/2023 14:42'
$ws.Range('T144').Value = 4.65
$ws.Range('U144').Value = '21/10/2023 15:56'
$ws.Range('V144').Value = 'https://www.betexplorer.com/football/england/championship/huddersfield-qpr/nBchU1vE/'

# Row 145 <- source data from former row 146
$ws.Range('F145').Value = 'Huddersfield'
$ws.Range('G145').Value = 0
$ws.Range('H145').Value = 'Cardiff'
$ws.Range('I145').Value = 4
$ws.Range('J145').Value = 2.37
$ws.Range('K145').Value = '20/10/2023 16:12'
$ws.Range('L145').Value = 2.52
$ws.Range('M145').Value = '24/10/2023 20:19'
$ws.Range('N145').Value = 3.41
$ws.Range('O145').Value = '20/10/2023 16:12'
$ws.Range('P145').Value = 3.17
$ws.Range('Q145').Value = '24/10/2023 20:41'
$ws.Range('R145').Value = 3.15
$ws.Range('S145').Value = '20/10/2023 16:12'
$ws.Range('T145').Value = 3.17
$ws.Range('U145').Value = '24/10/2023 20:19'
$ws.Range('V145').Value = 'https://www.betexplorer.com/football/england/championship/huddersfield-cardiff/vP8M9kFK/'

# Row 146 <- source data from former row 147
$ws.Range('F146').Value = 'Millwall'
$ws.Range('G146').Value = 1
$ws.Range('H146').Value = 'Blackburn'
$ws.Range('I146').Value = 2
$ws.Range('J146').Value = 2.26
$ws.Range('K146').Value = '20/10/2023 16:12'
$ws.Range('L146').Value = 2.68
$ws.Range('M146').Value = '24/10/2023 20:36'
$ws.Range('N146').Value = 3.57
$ws.Range('O146').Value = '20/10/2023 16:12'
$ws.Range('P146').Value = 3.49
$ws.Range('Q146').Value = '24/10/2023 20:28'
$ws.Range('R146').Value = 3.25
$ws.Range('S146').Value = '20/10/2023 16:12'
$ws.Range('T146').Value = 2.72
$ws.Range('U146').Value = '24/10/2023 20:36'
$ws.Range('V146').Value = 'https://www.betexplorer.com/football/england/championship/millwall-blackburn/2m7zsuOJ/'

# Row 147 <- source data from former row 148
$ws.Range('F147').Value = 'Norwich'
$ws.Range('G147').Value = 1
$ws.Range('H147').Value = 'Middlesbrough'
$ws.Range('I147').Value = 2
$ws.Range('J147').Value = 2.31
$ws.Range('K147').Value = '20/10/2023 16:12'
$ws.Range('L147').Value = 2.69
$ws.Range('M147').Value = '24/10/2023 20:43'
$ws.Range('N147').Value = 3.61
$ws.Range('O147').Value = '20/10/2023 16:12'
$ws.Range('P147').Value = 3.69
$ws.Range('Q147').Value = '24/10/2023 20:16'
$ws.Range('R147').Value = 3.12
$ws.Range('S147').Value = '20/10/2023 16:12'
$ws.Range('T147').Value = 2.6
$ws.Range('U147').Value = '24/10/2023 20:41'
$ws.Range('V147').Value = 'https://www.betexplorer.com/football/england/championship/norwich-middlesbrough/GjBvtawQ/'

# Row 148 <- source data from former row 145
$ws.Range('F148').Value = 'Swansea'
$ws.Range('G148').Value = 0
$ws.Range('H148').Value = 'Watford'
$ws.Range('I148').Value = 1
$ws.Range('J148').Value = 2.24
$ws.Range('K148').Value = '20/10/2023 16:12'
$ws.Range('L148').Value = 2.51
$ws.Range('M148').Value = '24/10/2023 20:31'
$ws.Range('N148').Value = 3.58
$ws.Range('O148').Value = '20/10/2023 16:12'
$ws.Range('P148').Value = 3.31
$ws.Range('Q148').Value = '24/10/2023 20:31'
$ws.Range('R148').Value = 3.27
$ws.Range('S148').Value = '20/10/2023 16:12'
$ws.Range('T148').Value = 3.06
$ws.Range('U148').Value = '24/10/2023 20:27'
$ws.Range('V148').Value = 'https://www.betexplorer.com/football/england/championship/swansea-watford/IF7Q89UQ/'

# --- Step 2: Append new rows 151-156 ---
$ws.Range('A150:V150').Copy($ws.Range('A151:V156'))

# New row 151
$ws.Range('A151').Value = 150
$ws.Range('B151').Value = 'england'
$ws.Range('C151').Value = 'championship'
$ws.Range('D151').Value = '2023-2024'
$ws.Range('E151').Value = 45224.86458333334
$ws.Range('F151').Value = 'Rotherham'
$ws.Range('G151').Value = 2
$ws.Range('H151').Value = 'Coventry'
$ws.Range('I151').Value = 0
$ws.Range('J151').Value = 3.64
$ws.Range('K151').Value = '20/10/2023 16:12'
$ws.Range('L151').Value = 5.34
$ws.Range('M151').Value = '25/10/2023 20:43'
$ws.Range('N151').Value = 3.66
$ws.Range('O151').Value = '20/10/2023 16:12'
$ws.Range('P151').Value = 3.86
$ws.Range('Q151').Value = '25/10/2023 20:43'
$ws.Range('R151').Value = 2.06
$ws.Range('S151').Value = '20/10/2023 16:12'
$ws.Range('T151').Value = 1.71
$ws.Range('U151').Value = '25/10/2023 20:43'
$ws.Range('V151').Value = 'https://www.betexplorer.com/football/england/championship/rotherham-coventry/EogRWf0C/'

# New row 152
$ws.Range('A152').Value = 151
$ws.Range('B152').Value = 'england'
$ws.Range('C152').Value = 'championship'
$ws.Range('D152').Value = '2023-2024'
$ws.Range('E152').Value = 45224.86458333334
$ws.Range('F152').Value = 'Birmingham'
$ws.Range('G152').Value = 0
$ws.Range('H152').Value = 'Hull'
$ws.Range('I152').Value = 2
$ws.Range('J152').Value = 2.16
$ws.Range('K152').Value = '20/10/2023 16:12'
$ws.Range('L152').Value = 2.53
$ws.Range('M152').Value = '25/10/2023 20:41'
$ws.Range('N152').Value = 3.5
$ws.Range('O152').Value = '20/10/2023 16:12'
$ws.Range('P152').Value = 3.34
$ws.Range('Q152').Value = '25/10/2023 20:37'
$ws.Range('R152').Value = 3.55
$ws.Range('S152').Value = '20/10/2023 16:12'
$ws.Range('T152').Value = 3.01
$ws.Range('U152').Value = '25/10/2023 20:37'
$ws.Range('V152').Value = 'https://www.betexplorer.com/football/england/championship/birmingham-hull-city/CMdAzw1m/'

# New row 153
$ws.Range('A153').Value = 152
$ws.Range('B153').Value = 'england'
$ws.Range('C153').Value = 'championship'
$ws.Range('D153').Value = '2023-2024'
$ws.Range('E153').Value = 45224.86458333334
$ws.Range('F153').Value = 'Bristol City'
$ws.Range('G153').Value = 0
$ws.Range('H153').Value = 'Ipswich'
$ws.Range('I153').Value = 1
$ws.Range('J153').Value = 3.03
$ws.Range('K153').Value = '20/10/2023 16:12'
$ws.Range('L153').Value = 4.03
$ws.Range('M153').Value = '25/10/2023 20:41'
$ws.Range('N153').Value = 3.57
$ws.Range('O153').Value = '20/10/2023 16:12'
$ws.Range('P153').Value = 3.87
$ws.Range('Q153').Value = '25/10/2023 20:39'
$ws.Range('R153').Value = 2.38
$ws.Range('S153').Value = '20/10/2023 16:12'
$ws.Range('T153').Value = 1.91
$ws.Range('U153').Value = '25/10/2023 20:41'
$ws.Range('V153').Value = 'https://www.betexplorer.com/football/england/championship/bristol-city-ipswich/UciFZdGg/'

# New row 154
$ws.Range('A154').Value = 153
$ws.Range('B154').Value = 'england'
$ws.Range('C154').Value = 'championship'
$ws.Range('D154').Value = '2023-2024'
$ws.Range('E154').Value = 45224.86458333334
$ws.Range('F154').Value = 'Plymouth'
$ws.Range('G154').Value = 3
$ws.Range('H154').Value = 'Sheffield Wed'
$ws.Range('I154').Value = 0
$ws.Range('J154').Value = 2.05
$ws.Range('K154').Value = '20/10/2023 16:12'
$ws.Range('L154').Value = 2.35
$ws.Range('M154').Value = '25/10/2023 20:43'
$ws.Range('N154').Value = 3.72
$ws.Range('O154').Value = '20/10/2023 16:12'
$ws.Range('P154').Value = 3.49
$ws.Range('Q154').Value = '25/10/2023 20:43'
$ws.Range('R154').Value = 3.64
$ws.Range('S154').Value = '20/10/2023 16:12'
$ws.Range('T154').Value = 3.17
$ws.Range('U154').Value = '25/10/2023 20:43'
$ws.Range('V154').Value = 'https://www.betexplorer.com/football/england/championship/plymouth-sheffield-wed/8ljJYGVa/'

# New row 155
$ws.Range('A155').Value = 154
$ws.Range('B155').Value = 'england'
$ws.Range('C155').Value = 'championship'
$ws.Range('D155').Value = '2023-2024'
$ws.Range('E155').Value = 45224.86458333334
$ws.Range('F155').Value = 'Preston'
$ws.Range('G155').Value = 2
$ws.Range('H155').Value = 'Southampton'
$ws.Range('I155').Value = 2
$ws.Range('J155').Value = 3.01
$ws.Range('K155').Value = '20/10/2023 16:12'
$ws.Range('L155').Value = 3.18
$ws.Range('M155').Value = '25/10/2023 20:43'
$ws.Range('N155').Value = 3.52
$ws.Range('O155').Value = '20/10/2023 16:12'
$ws.Range('P155').Value = 3.44
$ws.Range('Q155').Value = '25/10/2023 20:42'
$ws.Range('R155').Value = 2.42
$ws.Range('S155').Value = '20/10/2023 16:12'
$ws.Range('T155').Value = 2.37
$ws.Range('U155').Value = '25/10/2023 20:43'
$ws.Range('V155').Value = 'https://www.betexplorer.com/football/england/championship/preston-southampton/2BkNXzo6/'

# New row 156
$ws.Range('A156').Value = 155
$ws.Range('B156').Value = 'england'
$ws.Range('C156').Value = 'championship'
$ws.Range('D156').Value = '2023-2024'
$ws.Range('E156').Value = 45224.875
$ws.Range('F156').Value = 'Stoke'
$ws.Range('G156').Value = 1
$ws.Range('H156').Value = 'Leeds'
$ws.Range('I156').Value = 0
$ws.Range('J156').Value = 3.28
$ws.Range('K156').Value = '20/10/2023 16:12'
$ws.Range('L156').Value = 4.19
$ws.Range('M156').Value = '25/10/2023 20:15'
$ws.Range('N156').Value = 3.61
$ws.Range('O156').Value = '20/10/2023 16:12'
$ws.Range('P156').Value = 3.83
$ws.Range('Q156').Value = '25/10/2023 20:46'
$ws.Range('R156').Value = 2.21
$ws.Range('S156').Value = '20/10/2023 16:12'
$ws.Range('T156').Value = 1.88
$ws.Range('U156').Value = '25/10/2023 20:15'
$ws.Range('V156').Value = 'https://www.betexplorer.com/football/england/championship/stoke-city-leeds/dArWVEFI/'
